$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $a = $val -split ",\s*"
        $n = $a.Length
        if ($n -gt 1) {
            for ($i = 0; $i -lt $n; $i++) {
                for ($j = 0; $j -lt $n - $i - 1; $j++) {
                    $cax = $a[$j].ToCharArray()
                    $cay = $a[$j+1].ToCharArray()
                    $minLen = [Math]::Min($cax.Length, $cay.Length)
                    $diff = 0
                    for ($k = 0; $k -lt $minLen; $k++) {
                        $d = [int]$cax[$k] - [int]$cay[$k]
                        if ($d -ne 0 -and $diff -eq 0) { $diff = $d }
                    }
                    if ($diff -eq 0) { $diff = $cax.Length - $cay.Length }
                    if ($diff -lt 0) {
                        $tmp = $a[$j]
                        $a[$j] = $a[$j+1]
                        $a[$j+1] = $tmp
                    }
                }
            }
            $cell.Value2 = [string]::Join(", ", $a)
        }
    }
}
